$wb = $excel.ActiveWorkbook

# --- Tasks sheet (sheet1): add a new "Cost" column (D) ---
$tasks = $wb.Worksheets.Item("Tasks")

$tasks.Range("D1").Value = "Cost"
for ($r = 2; $r -le 23; $r++) {
    $tasks.Cells.Item($r, 4).Value = 1
}

# --- Update the saved view/selection state for each sheet ---
$resources = $wb.Worksheets.Item("Resources")
$resources.Range("E33").Select() | Out-Null

$tasks.Activate() | Out-Null
$tasks.Range("D2").Select() | Out-Null
